$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 38194725
$ws.Range("B2").Value = 18852567
$ws.Range("C2").Value = 54607854
$ws.Range("D2").Value = 31865985
$ws.Range("E2").Value = 16413129
$ws.Range("F2").Value = 42.97
$ws.Range("G2").Value = 13013418
$ws.Range("H2").Value = 69.03
